# Applies the Linea 141 scrape refresh (commit: "Horarios actualizados Linea 141 - 1159")
# Data snapshot timestamp moves from 11:54:47 -> 12:26:44 across all three sheets.
$wb = $excel.ActiveWorkbook

# ================= Sheet 1: LP1912 =================
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 12:26:44"
$ws1.Cells.Item(3, 1).Value = "Total filas: 176"

# A handful of earlier rows had their "Linea" (column C) values swapped between
# two adjacent rows that share the same Hora_Scrap/Hora_Llegada/Minutos.
$ws1.Cells.Item(22, 3).Value = "14_ABASTO"
$ws1.Cells.Item(23, 3).Value = "215C_EL PATO"

$ws1.Cells.Item(69, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(70, 3).Value = "84_COLONIA URQUIZA-ESC 49"

$ws1.Cells.Item(76, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(77, 3).Value = "16_SANTA ANA"

$ws1.Cells.Item(82, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(83, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(84, 3).Value = "17_ROMERO"

$ws1.Cells.Item(118, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(119, 3).Value = "10_OLMOS"

$ws1.Cells.Item(142, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(143, 3).Value = "23_HERNANDEZ"

$ws1.Cells.Item(149, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(150, 3).Value = "14_ABASTO"

# Rows 152-167 are rewritten with refreshed scrape data, and rows 168-181 are newly
# appended (total rows grows from 162 to 176, dimension A1:E167 -> A1:E181).
$ws1.Cells.Item(152, 1).Value = "12:26:44"
$ws1.Cells.Item(152, 2).Value = "12:27"
$ws1.Cells.Item(152, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(152, 4).Value = 1
$ws1.Cells.Item(152, 5).Value = "LP1912"

$ws1.Cells.Item(153, 1).Value = "12:26:44"
$ws1.Cells.Item(153, 2).Value = "12:34"
$ws1.Cells.Item(153, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(153, 4).Value = 8
$ws1.Cells.Item(153, 5).Value = "LP1912"

$ws1.Cells.Item(154, 1).Value = "12:26:44"
$ws1.Cells.Item(154, 2).Value = "12:34"
$ws1.Cells.Item(154, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(154, 4).Value = 8
$ws1.Cells.Item(154, 5).Value = "LP1912"

$ws1.Cells.Item(155, 1).Value = "12:26:44"
$ws1.Cells.Item(155, 2).Value = "12:36"
$ws1.Cells.Item(155, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(155, 4).Value = 10
$ws1.Cells.Item(155, 5).Value = "LP1912"

$ws1.Cells.Item(156, 1).Value = "12:26:44"
$ws1.Cells.Item(156, 2).Value = "12:36"
$ws1.Cells.Item(156, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(156, 4).Value = 10
$ws1.Cells.Item(156, 5).Value = "LP1912"

$ws1.Cells.Item(157, 1).Value = "12:26:44"
$ws1.Cells.Item(157, 2).Value = "12:38"
$ws1.Cells.Item(157, 3).Value = "17_179 Y 38"
$ws1.Cells.Item(157, 4).Value = 12
$ws1.Cells.Item(157, 5).Value = "LP1912"

$ws1.Cells.Item(158, 1).Value = "12:26:44"
$ws1.Cells.Item(158, 2).Value = "12:41"
$ws1.Cells.Item(158, 3).Value = "10_OLMOS"
$ws1.Cells.Item(158, 4).Value = 15
$ws1.Cells.Item(158, 5).Value = "LP1912"

$ws1.Cells.Item(159, 1).Value = "12:26:44"
$ws1.Cells.Item(159, 2).Value = "12:48"
$ws1.Cells.Item(159, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(159, 4).Value = 22
$ws1.Cells.Item(159, 5).Value = "LP1912"

$ws1.Cells.Item(160, 1).Value = "12:26:44"
$ws1.Cells.Item(160, 2).Value = "12:50"
$ws1.Cells.Item(160, 3).Value = "15_ABASTO"
$ws1.Cells.Item(160, 4).Value = 24
$ws1.Cells.Item(160, 5).Value = "LP1912"

$ws1.Cells.Item(161, 1).Value = "12:26:44"
$ws1.Cells.Item(161, 2).Value = "12:55"
$ws1.Cells.Item(161, 3).Value = "10_OLMOS"
$ws1.Cells.Item(161, 4).Value = 29
$ws1.Cells.Item(161, 5).Value = "LP1912"

$ws1.Cells.Item(162, 1).Value = "12:26:44"
$ws1.Cells.Item(162, 2).Value = "13:02"
$ws1.Cells.Item(162, 3).Value = "15_ABASTO"
$ws1.Cells.Item(162, 4).Value = 36
$ws1.Cells.Item(162, 5).Value = "LP1912"

$ws1.Cells.Item(163, 1).Value = "12:26:44"
$ws1.Cells.Item(163, 2).Value = "13:02"
$ws1.Cells.Item(163, 3).Value = "14_ABASTO"
$ws1.Cells.Item(163, 4).Value = 36
$ws1.Cells.Item(163, 5).Value = "LP1912"

$ws1.Cells.Item(164, 1).Value = "12:26:44"
$ws1.Cells.Item(164, 2).Value = "13:06"
$ws1.Cells.Item(164, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(164, 4).Value = 40
$ws1.Cells.Item(164, 5).Value = "LP1912"

$ws1.Cells.Item(165, 1).Value = "12:26:44"
$ws1.Cells.Item(165, 2).Value = "13:08"
$ws1.Cells.Item(165, 3).Value = "10_OLMOS"
$ws1.Cells.Item(165, 4).Value = 42
$ws1.Cells.Item(165, 5).Value = "LP1912"

$ws1.Cells.Item(166, 1).Value = "12:26:44"
$ws1.Cells.Item(166, 2).Value = "13:13"
$ws1.Cells.Item(166, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(166, 4).Value = 47
$ws1.Cells.Item(166, 5).Value = "LP1912"

$ws1.Cells.Item(167, 1).Value = "12:26:44"
$ws1.Cells.Item(167, 2).Value = "13:19"
$ws1.Cells.Item(167, 3).Value = "10_OLMOS"
$ws1.Cells.Item(167, 4).Value = 53
$ws1.Cells.Item(167, 5).Value = "LP1912"

$ws1.Cells.Item(168, 1).Value = "11:54:47"
$ws1.Cells.Item(168, 2).Value = "13:20"
$ws1.Cells.Item(168, 3).Value = "10_OLMOS"
$ws1.Cells.Item(168, 4).Value = 86
$ws1.Cells.Item(168, 5).Value = "LP1912"

$ws1.Cells.Item(169, 1).Value = "12:26:44"
$ws1.Cells.Item(169, 2).Value = "13:21"
$ws1.Cells.Item(169, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(169, 4).Value = 55
$ws1.Cells.Item(169, 5).Value = "LP1912"

$ws1.Cells.Item(170, 1).Value = "12:26:44"
$ws1.Cells.Item(170, 2).Value = "13:26"
$ws1.Cells.Item(170, 3).Value = "14_ABASTO"
$ws1.Cells.Item(170, 4).Value = 60
$ws1.Cells.Item(170, 5).Value = "LP1912"

$ws1.Cells.Item(171, 1).Value = "11:54:47"
$ws1.Cells.Item(171, 2).Value = "13:26"
$ws1.Cells.Item(171, 3).Value = "15_ABASTO"
$ws1.Cells.Item(171, 4).Value = 92
$ws1.Cells.Item(171, 5).Value = "LP1912"

$ws1.Cells.Item(172, 1).Value = "12:26:44"
$ws1.Cells.Item(172, 2).Value = "13:34"
$ws1.Cells.Item(172, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(172, 4).Value = 68
$ws1.Cells.Item(172, 5).Value = "LP1912"

$ws1.Cells.Item(173, 1).Value = "11:54:47"
$ws1.Cells.Item(173, 2).Value = "13:36"
$ws1.Cells.Item(173, 3).Value = "15_ABASTO"
$ws1.Cells.Item(173, 4).Value = 102
$ws1.Cells.Item(173, 5).Value = "LP1912"

$ws1.Cells.Item(174, 1).Value = "12:26:44"
$ws1.Cells.Item(174, 2).Value = "13:46"
$ws1.Cells.Item(174, 3).Value = "17_ROMERO"
$ws1.Cells.Item(174, 4).Value = 80
$ws1.Cells.Item(174, 5).Value = "LP1912"

$ws1.Cells.Item(175, 1).Value = "12:26:44"
$ws1.Cells.Item(175, 2).Value = "13:50"
$ws1.Cells.Item(175, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(175, 4).Value = 84
$ws1.Cells.Item(175, 5).Value = "LP1912"

$ws1.Cells.Item(176, 1).Value = "12:26:44"
$ws1.Cells.Item(176, 2).Value = "13:56"
$ws1.Cells.Item(176, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(176, 4).Value = 90
$ws1.Cells.Item(176, 5).Value = "LP1912"

$ws1.Cells.Item(177, 1).Value = "12:26:44"
$ws1.Cells.Item(177, 2).Value = "13:56"
$ws1.Cells.Item(177, 3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(177, 4).Value = 90
$ws1.Cells.Item(177, 5).Value = "LP1912"

$ws1.Cells.Item(178, 1).Value = "12:26:44"
$ws1.Cells.Item(178, 2).Value = "14:04"
$ws1.Cells.Item(178, 3).Value = "17_ROMERO"
$ws1.Cells.Item(178, 4).Value = 98
$ws1.Cells.Item(178, 5).Value = "LP1912"

$ws1.Cells.Item(179, 1).Value = "12:26:44"
$ws1.Cells.Item(179, 2).Value = "14:16"
$ws1.Cells.Item(179, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(179, 4).Value = 110
$ws1.Cells.Item(179, 5).Value = "LP1912"

$ws1.Cells.Item(180, 1).Value = "12:26:44"
$ws1.Cells.Item(180, 2).Value = "14:19"
$ws1.Cells.Item(180, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(180, 4).Value = 113
$ws1.Cells.Item(180, 5).Value = "LP1912"

$ws1.Cells.Item(181, 1).Value = "12:26:44"
$ws1.Cells.Item(181, 2).Value = "14:21"
$ws1.Cells.Item(181, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(181, 4).Value = 115
$ws1.Cells.Item(181, 5).Value = "LP1912"

# ================= Sheet 2: LP1912-215 =================
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 12:26:44"
$ws2.Cells.Item(3, 1).Value = "Total filas: 19"

$ws2.Cells.Item(22, 1).Value = "12:26:44"
$ws2.Cells.Item(22, 4).Value = 47

$ws2.Cells.Item(23, 1).Value = "12:26:44"
$ws2.Cells.Item(23, 4).Value = 84

# New row 24 appended (dimension A1:E23 -> A1:E24)
$ws2.Cells.Item(24, 1).Value = "12:26:44"
$ws2.Cells.Item(24, 2).Value = "14:19"
$ws2.Cells.Item(24, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(24, 4).Value = 113
$ws2.Cells.Item(24, 5).Value = "LP1912"

# ================= Sheet 3: 6203-6173 =================
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 12:26:44"
$ws3.Cells.Item(3, 1).Value = "Total filas: 26"

$ws3.Cells.Item(29, 1).Value = "12:26:44"
$ws3.Cells.Item(29, 4).Value = 27

$ws3.Cells.Item(30, 1).Value = "12:26:44"
$ws3.Cells.Item(30, 4).Value = 64

# New row 31 appended (dimension A1:E30 -> A1:E31)
$ws3.Cells.Item(31, 1).Value = "12:26:44"
$ws3.Cells.Item(31, 2).Value = "14:09"
$ws3.Cells.Item(31, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(31, 4).Value = 103
$ws3.Cells.Item(31, 5).Value = "L6173"

